$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 728; existing rows 728-789 shift down to 729-790.
$ws.Rows.Item(728).Insert()

# Populate the new row 728 with the new weekly record.
$ws.Range("A728").Value = 3
$ws.Range("B728").Value = "Femacal de La Calera"
$ws.Range("C728").Value = "Coquimbo"
$ws.Range("D728").Value = 45106
$ws.Range("E728").Value = 5
$ws.Range("F728").Value = 100112032
$ws.Range("G728").Value = "Zapallo italiano"
$ws.Range("H728").Value = "Sin especificar"
$ws.Range("I728").Value = "Primera"
$ws.Range("J728").Value = 250
$ws.Range("K728").Value = 11000
$ws.Range("L728").Value = 11500
$ws.Range("M728").Value = 11268
$ws.Range("N728").Value = "$/caja 60 unidades"
$ws.Range("O728").Value = "Región de Arica y Parinacota"
$ws.Range("P728").Value = 188
$ws.Range("Q728").Value = 60
$ws.Range("R728").Value = "Hortaliza"

# Match the date-number-format style used by the other rows in column D.
$ws.Range("D728").NumberFormat = $ws.Range("D729").NumberFormat
